$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 237, shifting existing rows 237:269 down to 238:270
$ws.Rows(237).Insert()

# Fill in the new row 237 with the inserted weekly record.
# Columns A,B,C,E,F,G,H,Q,R carry the same constant values as the surrounding rows.
$ws.Cells.Item(237, 1).Value = 3
$ws.Cells.Item(237, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(237, 3).Value = "Coquimbo"
$ws.Cells.Item(237, 4).Value = 44505
$ws.Cells.Item(237, 5).Value = 5
$ws.Cells.Item(237, 6).Value = 100112028
$ws.Cells.Item(237, 7).Value = "Sandia"
$ws.Cells.Item(237, 8).Value = "Sin especificar"
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 160
$ws.Cells.Item(237, 11).Value = 600
$ws.Cells.Item(237, 12).Value = 600
$ws.Cells.Item(237, 13).Value = 600
$ws.Cells.Item(237, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(237, 15).Value = "Perú"
$ws.Cells.Item(237, 16).Value = 600
$ws.Cells.Item(237, 17).Value = 1
$ws.Cells.Item(237, 18).Value = "Hortaliza"
